$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$labels = @(
    "BLDAS",
    "MALE",
    "FCRL6",
    "GDNF",
    "KRT19",
    "HAOX1",
    "REN",
    "MMP-1",
    "ARNT",
    "IL13",
    "TRAIL-R2",
    "MCP-2",
    "CXCL1",
    "PRKCQ",
    "TNFSF13B",
    "RARRES2",
    "SPON1",
    "hOSCAR",
    "DPP10"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
